# Swap the order of "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
# in the "Recorded By" column (column G) of the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "dnasr281@gmail.com, System"
$newText = "System, dnasr281@gmail.com"

# Use Excel's native Replace so only cells containing the exact text are
# touched (avoids materializing already-empty cells when iterating manually).
$col = $ws.Columns.Item(7)  # Column G = "Recorded By"
$col.Replace($oldText, $newText)
